# Add settler/settling production data (Residence building, levels 1-10)
# to the "fixed data" sheet, rows 27-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("fixed data")
$ws.Activate()

# Columns: B=lvl, C=Lumber, D=Clay, E=Iron, F=Crop, G=pop, H=CP,
# I=Time (h:mm:ss), J=Prod. increase. Column A is the building name
# ("Residence", shared string already present in the workbook).
$rows = @(
    @(1,  580,  460,  350,  180, 1, 2,  2000/86400,  2),
    @(2,  740,  590,  450,  230, 1, 3,  2620/86400,  8),
    @(3,  950,  755,  575,  295, 1, 3,  3340/86400,  18),
    @(4,  1215, 965,  735,  375, 1, 4,  4170/86400,  32),
    @(5,  1555, 1235, 940,  485, 1, 5,  5140/86400,  50),
    @(6,  1995, 1580, 1205, 620, 1, 6,  6260/86400,  72),
    @(7,  2550, 2025, 1540, 790, 1, 7,  7570/86400,  98),
    @(8,  3265, 2590, 1970, 1015,1, 9,  9080/86400,  128),
    @(9,  4180, 3315, 2520, 1295,1, 10, 10830/86400, 162),
    @(10, 5350, 4245, 3230, 1660,1, 12, 12860/86400, 200)
)

$r = 27
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = "Residence"
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 9).NumberFormat = "h:mm:ss"
    $ws.Cells.Item($r, 10).Value = $row[8]
    $r = $r + 1
}

# Match the saved view state: scrolled down one row, active cell N10.
$ws.Range("N10").Select()
